# Update the "取得日時" (retrieved datetime) timestamp in column A
# for the data rows on the "ランサーズ" sheet from 2025-12-23 06:30:49
# to 2025-12-23 06:39:18 (new append run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-23 06:39:18"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
